# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (quarterly fund-holdings detail, same
# shape as the existing 2021-Qx sheets) right before the "总计" summary sheet,
# then updates "总计" with a new leading row for 2022-Q1 (existing rows shift
# down by one).
#
# NOTE: worksheet references captured *before* Worksheets.Add() go stale once
# the sheet collection is mutated (they keep tracking the old positional
# index, not the sheet they were bound to) - so every sheet handle we keep
# using after the Add() call below is re-fetched by name afterwards.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "2022-Q1" sheet right before "总计" -----------------
$totalSheetForPositioning = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetForPositioning)
$newSheet.Name = "2022-Q1"

# Re-fetch every sheet handle we still need now that the collection changed.
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# --- 2. Header row ----------------------------------------------------------
$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

# --- 3. Fund rows -------------------------------------------------------------
# row 2: 010347
$newSheet.Cells.Item(2, 1).Value2 = 0
$newSheet.Cells.Item(2, 2).Value2 = "'010347"
$newSheet.Cells.Item(2, 3).Value2 = "农银汇理策略收益一年持有期混合"
$newSheet.Cells.Item(2, 4).Value2 = "'56.47"
$newSheet.Cells.Item(2, 5).Value2 = "'74.89"
$newSheet.Cells.Item(2, 6).Value2 = "'5.41"
$newSheet.Cells.Item(2, 7).Value2 = "'3.0550"
$newSheet.Cells.Item(2, 8).Value2 = 5

# row 3: 660010
$newSheet.Cells.Item(3, 1).Value2 = 1
$newSheet.Cells.Item(3, 2).Value2 = "'660010"
$newSheet.Cells.Item(3, 3).Value2 = "农银策略精选混合"
$newSheet.Cells.Item(3, 4).Value2 = "'31.49"
$newSheet.Cells.Item(3, 5).Value2 = "'75.60"
$newSheet.Cells.Item(3, 6).Value2 = "'5.29"
$newSheet.Cells.Item(3, 7).Value2 = "'1.6658"
$newSheet.Cells.Item(3, 8).Value2 = 5

# row 4: 000127
$newSheet.Cells.Item(4, 1).Value2 = 2
$newSheet.Cells.Item(4, 2).Value2 = "'000127"
$newSheet.Cells.Item(4, 3).Value2 = "农银行业领先混合"
$newSheet.Cells.Item(4, 4).Value2 = "'12.28"
$newSheet.Cells.Item(4, 5).Value2 = "'75.75"
$newSheet.Cells.Item(4, 6).Value2 = "'5.31"
$newSheet.Cells.Item(4, 7).Value2 = "'0.6521"
$newSheet.Cells.Item(4, 8).Value2 = 5

# row 5: 011817
$newSheet.Cells.Item(5, 1).Value2 = 3
$newSheet.Cells.Item(5, 2).Value2 = "'011817"
$newSheet.Cells.Item(5, 3).Value2 = "银华阿尔法混合型证券投资基金"
$newSheet.Cells.Item(5, 4).Value2 = "'10.82"
$newSheet.Cells.Item(5, 5).Value2 = "'67.94"
$newSheet.Cells.Item(5, 6).Value2 = "'4.35"
$newSheet.Cells.Item(5, 7).Value2 = "'0.4707"
$newSheet.Cells.Item(5, 8).Value2 = 5

# row 6: 501075
$newSheet.Cells.Item(6, 1).Value2 = 4
$newSheet.Cells.Item(6, 2).Value2 = "'501075"
$newSheet.Cells.Item(6, 3).Value2 = "万家科创主题3年封闭运作灵活配置混合A"
$newSheet.Cells.Item(6, 4).Value2 = "'17.95"
$newSheet.Cells.Item(6, 5).Value2 = "'79.41"
$newSheet.Cells.Item(6, 6).Value2 = "'2.30"
$newSheet.Cells.Item(6, 7).Value2 = "'0.4128"
$newSheet.Cells.Item(6, 8).Value2 = 10

# row 7: 008819
$newSheet.Cells.Item(7, 1).Value2 = 5
$newSheet.Cells.Item(7, 2).Value2 = "'008819"
$newSheet.Cells.Item(7, 3).Value2 = "农银汇理策略趋势混合"
$newSheet.Cells.Item(7, 4).Value2 = "'6.17"
$newSheet.Cells.Item(7, 5).Value2 = "'76.46"
$newSheet.Cells.Item(7, 6).Value2 = "'5.34"
$newSheet.Cells.Item(7, 7).Value2 = "'0.3295"
$newSheet.Cells.Item(7, 8).Value2 = 5

# row 8: 180001
$newSheet.Cells.Item(8, 1).Value2 = 6
$newSheet.Cells.Item(8, 2).Value2 = "'180001"
$newSheet.Cells.Item(8, 3).Value2 = "银华优势企业混合"
$newSheet.Cells.Item(8, 4).Value2 = "'6.66"
$newSheet.Cells.Item(8, 5).Value2 = "'67.80"
$newSheet.Cells.Item(8, 6).Value2 = "'4.44"
$newSheet.Cells.Item(8, 7).Value2 = "'0.2957"
$newSheet.Cells.Item(8, 8).Value2 = 6

# row 9: 090016
$newSheet.Cells.Item(9, 1).Value2 = 7
$newSheet.Cells.Item(9, 2).Value2 = "'090016"
$newSheet.Cells.Item(9, 3).Value2 = "大成消费主题混合"
$newSheet.Cells.Item(9, 4).Value2 = "'4.23"
$newSheet.Cells.Item(9, 5).Value2 = "'93.78"
$newSheet.Cells.Item(9, 6).Value2 = "'6.09"
$newSheet.Cells.Item(9, 7).Value2 = "'0.2576"
$newSheet.Cells.Item(9, 8).Value2 = 6

# row 10: 001163
$newSheet.Cells.Item(10, 1).Value2 = 8
$newSheet.Cells.Item(10, 2).Value2 = "'001163"
$newSheet.Cells.Item(10, 3).Value2 = "银华中国梦30股票"
$newSheet.Cells.Item(10, 4).Value2 = "'5.00"
$newSheet.Cells.Item(10, 5).Value2 = "'90.33"
$newSheet.Cells.Item(10, 6).Value2 = "'4.97"
$newSheet.Cells.Item(10, 7).Value2 = "'0.2485"
$newSheet.Cells.Item(10, 8).Value2 = 7

# row 11: 501070
$newSheet.Cells.Item(11, 1).Value2 = 9
$newSheet.Cells.Item(11, 2).Value2 = "'501070"
$newSheet.Cells.Item(11, 3).Value2 = "广发睿阳三年定期开放混合"
$newSheet.Cells.Item(11, 4).Value2 = "'7.06"
$newSheet.Cells.Item(11, 5).Value2 = "'50.14"
$newSheet.Cells.Item(11, 6).Value2 = "'2.85"
$newSheet.Cells.Item(11, 7).Value2 = "'0.2012"
$newSheet.Cells.Item(11, 8).Value2 = 7

# row 12: 002563
$newSheet.Cells.Item(12, 1).Value2 = 10
$newSheet.Cells.Item(12, 2).Value2 = "'002563"
$newSheet.Cells.Item(12, 3).Value2 = "泓德泓汇灵活配置混合"
$newSheet.Cells.Item(12, 4).Value2 = "'4.61"
$newSheet.Cells.Item(12, 5).Value2 = "'92.58"
$newSheet.Cells.Item(12, 6).Value2 = "'4.02"
$newSheet.Cells.Item(12, 7).Value2 = "'0.1853"
$newSheet.Cells.Item(12, 8).Value2 = 6

# row 13: 001695
$newSheet.Cells.Item(13, 1).Value2 = 11
$newSheet.Cells.Item(13, 2).Value2 = "'001695"
$newSheet.Cells.Item(13, 3).Value2 = "泓德泓业灵活配置混合"
$newSheet.Cells.Item(13, 4).Value2 = "'1.16"
$newSheet.Cells.Item(13, 5).Value2 = "'91.88"
$newSheet.Cells.Item(13, 6).Value2 = "'4.02"
$newSheet.Cells.Item(13, 7).Value2 = "'0.0466"
$newSheet.Cells.Item(13, 8).Value2 = 5

# row 14: 519093
$newSheet.Cells.Item(14, 1).Value2 = 12
$newSheet.Cells.Item(14, 2).Value2 = "'519093"
$newSheet.Cells.Item(14, 3).Value2 = "新华钻石品质企业混合"
$newSheet.Cells.Item(14, 4).Value2 = "'1.36"
$newSheet.Cells.Item(14, 5).Value2 = "'92.58"
$newSheet.Cells.Item(14, 6).Value2 = "'3.06"
$newSheet.Cells.Item(14, 7).Value2 = "'0.0416"
$newSheet.Cells.Item(14, 8).Value2 = 10

# row 15: 007861
$newSheet.Cells.Item(15, 1).Value2 = 13
$newSheet.Cells.Item(15, 2).Value2 = "'007861"
$newSheet.Cells.Item(15, 3).Value2 = "金元顺安医疗健康混合型证券投资基金A"
$newSheet.Cells.Item(15, 4).Value2 = "'0.52"
$newSheet.Cells.Item(15, 5).Value2 = "'86.80"
$newSheet.Cells.Item(15, 6).Value2 = "'4.20"
$newSheet.Cells.Item(15, 7).Value2 = "'0.0218"
$newSheet.Cells.Item(15, 8).Value2 = 5

# row 16: 007862
$newSheet.Cells.Item(16, 1).Value2 = 14
$newSheet.Cells.Item(16, 2).Value2 = "'007862"
$newSheet.Cells.Item(16, 3).Value2 = "金元顺安医疗健康混合型证券投资基金C"
$newSheet.Cells.Item(16, 4).Value2 = "'0.09"
$newSheet.Cells.Item(16, 5).Value2 = "'86.80"
$newSheet.Cells.Item(16, 6).Value2 = "'4.20"
$newSheet.Cells.Item(16, 7).Value2 = "'0.0038"
$newSheet.Cells.Item(16, 8).Value2 = 5

# row 17: 007501
$newSheet.Cells.Item(17, 1).Value2 = 15
$newSheet.Cells.Item(17, 2).Value2 = "'007501"
$newSheet.Cells.Item(17, 3).Value2 = "万家科创主题3年封闭运作灵活配置混合C"
# (fund-scale cell intentionally left blank for this row, per source data)
$newSheet.Cells.Item(17, 5).Value2 = "'79.41"
$newSheet.Cells.Item(17, 6).Value2 = "'2.30"
$newSheet.Cells.Item(17, 7).Value2 = 0
$newSheet.Cells.Item(17, 8).Value2 = 10

# row 18: 003739
$newSheet.Cells.Item(18, 1).Value2 = 16
$newSheet.Cells.Item(18, 2).Value2 = "'003739"
$newSheet.Cells.Item(18, 3).Value2 = "新华鑫弘灵活配置混合"
$newSheet.Cells.Item(18, 4).Value2 = "'0.01"
$newSheet.Cells.Item(18, 5).Value2 = "'42.81"
$newSheet.Cells.Item(18, 6).Value2 = "'0.45"
$newSheet.Cells.Item(18, 7).Value2 = "'0.0000"
$newSheet.Cells.Item(18, 8).Value2 = 10

# --- 4. Match the header / index-column style used by the other quarterly --
#        sheets (e.g. "2021-Q4"): bold, bordered, centered.
$q4Sheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2:A18").PasteSpecial(-4122)

# --- 5. Update the "总计" sheet: push existing rows down one slot and -------
#        insert the new 2022-Q1 summary line at the top.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

$totalSheet.Cells.Item(5, 1).Value2 = 3
$totalSheet.Cells.Item(5, 2).Value2 = $totalSheet.Cells.Item(4, 2).Value2
$totalSheet.Cells.Item(5, 3).Value2 = $totalSheet.Cells.Item(4, 3).Value2
$totalSheet.Cells.Item(5, 4).Value2 = $totalSheet.Cells.Item(4, 4).Value2

$totalSheet.Cells.Item(4, 1).Value2 = 2
$totalSheet.Cells.Item(4, 2).Value2 = $totalSheet.Cells.Item(3, 2).Value2
$totalSheet.Cells.Item(4, 3).Value2 = $totalSheet.Cells.Item(3, 3).Value2
$totalSheet.Cells.Item(4, 4).Value2 = $totalSheet.Cells.Item(3, 4).Value2

$totalSheet.Cells.Item(3, 1).Value2 = 1
$totalSheet.Cells.Item(3, 2).Value2 = $totalSheet.Cells.Item(2, 2).Value2
$totalSheet.Cells.Item(3, 3).Value2 = $totalSheet.Cells.Item(2, 3).Value2
$totalSheet.Cells.Item(3, 4).Value2 = $totalSheet.Cells.Item(2, 4).Value2

$totalSheet.Cells.Item(2, 1).Value2 = 0
$totalSheet.Cells.Item(2, 2).Value2 = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value2 = 17
$totalSheet.Cells.Item(2, 4).Value2 = 7.89

